# add test in GASolver
# Bump a set of lookup-table totals in column D by 1 (test-fixture refresh)
# and move the active selection to L17, matching the authored commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows whose column-D value increments by 1.
$rowsToBump = 2,3,6,7,10,11,12,15,16,18,21,23,25,27,28,31,32,35,36,37,40,41,43,46,48,50

foreach ($r in $rowsToBump) {
    $cell = $ws.Cells.Item($r, 4)
    $cell.Value = $cell.Value2 + 1
}

# Move the active selection (was F1:F1048576 with active cell F1).
$ws.Range("L17").Select()
